$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (F2_Clientes_Compradores_Score) = 100 for every data row (2-13).
# Column I (F3_Pedidos_Por_Dia) mirrors column C's value for that row.
# Both newly-filled cells pick up the same "red font + border" look already
# used by the other numeric cells in the row (style index 3 in the original
# file), which Excel applies automatically once we set a matching font color.

$gValues = @(100, 100, 100, 100, 100, 100, 100, 100, 100, 100, 100, 100)
$iValues = @(100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111)

for ($row = 2; $row -le 13; $row++) {
    $idx = $row - 2

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $gValues[$idx]
    $gCell.Font.Color = 255

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value = $iValues[$idx]
    $iCell.Font.Color = 255
}

# Update the view: scroll right a bit and move the selection to H2.
$ws.Range("H2").Select()
